# B6-PowerPoint.pptx edit
#
# 1) Three tables (on slides 14, 15, 16) get their table style switched
#    from the deck-local default style {A9002C06-2DEA-4F1D-AE00-C42546E71154}
#    to the built-in style {2450B5FB-CA1B-48A2-97FB-D554BE541A74}.
#
# 2) The theme bound to the slide master (ppt/theme/theme1.xml, the
#    "Integral"/"Red Violet" color set) is swapped for the "Office Theme"
#    color set (the twelve theme colors that previously lived in
#    ppt/theme/theme2.xml / the notes-master theme).

$p = $ppt.ActivePresentation

# --- 1) Re-style the three tables -----------------------------------------
$newTableStyleId = "{2450B5FB-CA1B-48A2-97FB-D554BE541A74}"

foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    $tableShape = $slide.Shapes.Item(1)
    if ($tableShape.HasTable) {
        $tableShape.Table.ApplyStyle($newTableStyleId)
    }
}

# --- 2) Swap the slide-master theme colors for the Office Theme palette ---
# (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink — in that order, as
# VBA/COM-style BGR-packed RGB() integers.)
$officeThemeColors = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

$slideOne = $p.Slides.Item(1)
$themeColors = $slideOne.ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = $officeThemeColors[$i - 1]
}
